# "added some fake data to the data"
# Adds two new Q&A rows (rows 3 & 4) to the fake_data sheet, a third
# Answer column-group (Q:U) in the header row, and the formatting that
# goes along with the new data (a dark-grey font for one answer, a
# wrapped/tall cell for a long multi-line question, and portrait page
# setup).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row additions (Answer3 column group) ----
$ws.Range("Q1").Value = "Answer3"
$ws.Range("R1").Value = "Answer3_Username"
$ws.Range("S1").Value = "Answer3_Time"
$ws.Range("T1").Value = "Answer3_Upvotes"
$ws.Range("U1").Value = "Answer3_Downvotes"

# ---- Row 3: brussel sprouts / coconut oil question ----
$ws.Range("A3").Value = "I just love brussel sprouts. And this guy (http://mywholefoodlife.com/2014/07/14/smokey-roasted-brussel-sprouts/) is too! The recipe he uses calls for brussel sprouts, 1 large sweet potato, 2 tablespoons of coconut oil, s tablesppons of maple syrup, 1/2 teaspoon of chili powder, and 1/2 teaspoon of sea salt. ....The coconut oil really does it for me. It's just so fatty! (Whatevs about the maple syrup - YOLO.) What can I use instead of the 2 tablespoons of coconut oil???"
$ws.Range("B3").Value = "maplesyrup<333"
$ws.Range("C3").Value = "11/16/2014  13:29:45 PM"
$ws.Range("C3").NumberFormat = "m/d/yy h:mm"
$ws.Range("D3").Value = "Roasting"
$ws.Range("E3").Value = "None"
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = "Coconut oil has a lot of unhealthy saturated fats, while olive oil is a healthier alternative with more unsaturated fats."
$ws.Range("G3").Font.Color = 2236962
$ws.Range("H3").Value = "ilikehealthyfoods"
$ws.Range("I3").Value = "11/16/2014  14:05:51 PM"
$ws.Range("I3").NumberFormat = "m/d/yy h:mm"
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = "Use 2 tablespoons of canola oil! It's less expensive than olive oil."
$ws.Range("N3").Value = 41959.613449074073
$ws.Range("N3").NumberFormat = "m/d/yy h:mm"
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = "Nobody likes brussel sprouts. Make brownies."
$ws.Range("R3").Value = "downwithvegetables"
$ws.Range("S3").Value = 41959.628101851849
$ws.Range("S3").NumberFormat = "m/d/yy h:mm"
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 4

# ---- Row 4: banana bread / lactose-intolerance question ----
$ws.Range("A4").Value = "Bananas! Bananas! Don't you love banana bread?   4 ripe bananas,
1 cup butter,
2 cups sugar,
2 large eggs,
1/2 cup low-fat mayonaise,
1/2 cup milk,
1 teaspoon baking soda,
3 cups all-purpose flour,
1 teaspoon vanilla,
1/2 teaspoons salt. But I'm lactose-intolerant. What can I do about the milk and the butter? Thank you!"
$ws.Range("A4").WrapText = $true
$ws.Range("B4").Value = "ughMilkWhyDoYouHurtMeSo"
$ws.Range("C4").Value = "11/16/2014  14:29:45 PM"
$ws.Range("C4").NumberFormat = "m/d/yy h:mm"
$ws.Range("D4").Value = "Baking"
$ws.Range("E4").Value = "Lactose_Intolerant"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "Use margarine and soy milk instead of dairy products. I keep kosher and often can't cook using dairy. Enjoy :)"
$ws.Range("H4").Value = "separate_milk_and_eat"
$ws.Range("I4").Value = "11/16/2014  14:38:06 PM"
$ws.Range("I4").NumberFormat = "m/d/yy h:mm"
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 0

# Row 4 renders as a tall wrapped row in the source workbook.
$ws.Rows.Item(4).RowHeight = 195

# New column (S) picks up an explicit width, same as the other answer
# columns.
$ws.Columns("S").ColumnWidth = 15.45

# Sheet now prints in portrait orientation.
$ws.PageSetup.Orientation = 1
